$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B, shifting the existing
# "Jun_13" (B) and "Jun_10" (C) history columns right to D and E.
$ws.Range("B:C").Insert()

# New column headers (most recent dates first).
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the new columns with the "UN" (unchanged) status used elsewhere
# in the sheet, for every data row.
$ws.Range("B2:C27").Value = "UN"
